$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-03 -> 2023-10-04, i.e. serial 45202 -> 45203) for every data
# row from row 2 through row 264.
$lastRow = 264
$ws.Range("C2:C$lastRow").Value = 45203
